# Applies the weekly fruit/hortaliza price update for the Papaya sheet.
# The edit reshuffles the Fecha/Calidad/Volumen/Precio* data among rows 2-12
# (all other columns - Mercado, Region, Producto, etc. - stay identical).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ("after") values for columns D, L, M, N, O, P, S for each data row (2-12).
$rows = @{
    2  = @{ D = 44503; L = "Primera"; M = 60;  N = 30000; O = 30000; P = 30000; S = 3000 }
    3  = @{ D = 44503; L = "Segunda"; M = 50;  N = 25000; O = 25000; P = 25000; S = 2500 }
    4  = @{ D = 44517; L = "Especial"; M = 100; N = 27000; O = 27000; P = 27000; S = 2700 }
    5  = @{ D = 44517; L = "Primera"; M = 30;  N = 25000; O = 25000; P = 25000; S = 2500 }
    6  = @{ D = 44432; L = "Primera"; M = 20;  N = 20000; O = 20000; P = 20000; S = 2000 }
    7  = @{ D = 44476; L = "Primera"; M = 120; N = 20000; O = 20000; P = 20000; S = 2000 }
    8  = @{ D = 44434; L = "Primera"; M = 20;  N = 20000; O = 20000; P = 20000; S = 2000 }
    9  = @{ D = 44435; L = "Primera"; M = 40;  N = 20000; O = 20000; P = 20000; S = 2000 }
    10 = @{ D = 44466; L = "Primera"; M = 60;  N = 20000; O = 20000; P = 20000; S = 2000 }
    11 = @{ D = 44511; L = "Primera"; M = 120; N = 28000; O = 28000; P = 28000; S = 2800 }
    12 = @{ D = 44473; L = "Primera"; M = 180; N = 20000; O = 20000; P = 20000; S = 2000 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals.S   # S - Precio $/Kg
}
